# Edit supervision.xlsx per commit:
# "Added plos one to revs, youtube link for brazilian talk, and reorganised folders"
#
# Effective content change on the "supervision" sheet:
#  - The "Biology" / "Music Pedagogy" rows (previously rows 20-21, with their
#    two "why" rows 22 removed to after) are moved up to directly follow the
#    "MSc in Psychology" row (now rows 6-8).
#  - The "Psychology" / Undergraduate row (previously row 6) moves down to row 9.
#  - Four new supervised students for 2021-2022 are inserted right after row 9:
#      Angie Alejandra Lozano Sanjuan (2021 - 2022)
#      Daniela Martínez Franco (2021 - 2022)
#      Mariana Saavedra Botero (2021 - 2022)
#      John Jairo Rubio (2021 - 2022)
#  - Every other "why" row shifts down accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("supervision")

# Clear out the existing data area (keep header row 1 untouched) before
# rewriting it in the new row order.
$ws.Range("A2:E30").ClearContents()

# Row 2 - unchanged
$ws.Range("A2").Value = "PhD in Neuroscience"
$ws.Range("B2").Value = "2015 - 2018"
$ws.Range("C2").Value = "Postgraduate"
$ws.Range("D2").Value = "University of Valencia, Spain"
$ws.Range("E2").Value = "Milena Vásquez-Amézquita. Supervised together with  Alicia Salvador"

# Row 3 - unchanged
$ws.Range("A3").Value = "Professional Doctorate in Counselling Psychology"
$ws.Range("B3").Value = "2015 - 2018"
$ws.Range("D3").Value = "University of East London, UK"
$ws.Range("E3").Value = "Francisco Javier Flores. Supervised together with Lisa Chiara Fellin"

# Row 4 - unchanged
$ws.Range("A4").Value = "Psychological Research Methods (Evolutionary Psychology) MSc"
$ws.Range("B4").Value = "2013 - 2014"
$ws.Range("D4").Value = "University of Stirling, UK"
$ws.Range("E4").Value = "Julia Sanz-Vidania. Supervised together with S Craig Roberts"

# Row 5 - unchanged
$ws.Range("A5").Value = "MSc in Psychology"
$ws.Range("B5").Value = "2019 - 2020"
$ws.Range("D5").Value = "Universidad El Bosque, Colombia"
$ws.Range("E5").Value = "Adrián Acosta Guerrero. Supervised together with Milena Vásquez-Amézquita"

# Row 6 - moved up from old row 20 (Biology)
$ws.Range("A6").Value = "Biology"
$ws.Range("B6").Value = "2017 - 2018"
$ws.Range("D6").Value = "Universidad El Bosque, Colombia"
$ws.Range("E6").Value = "Maria Alejandra Abello Mozo  (2017 - 2018)"

# Row 7 - moved up from old row 21 (Music Pedagogy)
$ws.Range("A7").Value = "Music Pedagogy"
$ws.Range("B7").Value = "2017 - 2019"
$ws.Range("D7").Value = "Universidad Pedagógica Nacional, Colombia"
$ws.Range("E7").Value = "Natalia Elízabeth Moreno Buitrago (2017 ‑ 2019)"

# Row 8 - moved up from old row 22
$ws.Range("E8").Value = "Juan Felipe Pérez Ariza (2017 ‑ 2019)"

# Row 9 - moved down from old row 6 (Psychology / Undergraduate)
$ws.Range("A9").Value = "Psychology"
$ws.Range("B9").Value = "Since 2015"
$ws.Range("C9").Value = "Undergraduate"
$ws.Range("D9").Value = "Universidad El Bosque, Colombia"
$ws.Range("E9").Value = "Andrés Castellanos-Chacón (2017 - 2018; teaching supervision 2019 - Present)"

# Rows 10-13 - four new students added (2021 - 2022)
$ws.Range("E10").Value = "Angie Alejandra Lozano Sanjuan (2021 - 2022)"
$ws.Range("E11").Value = "Daniela Martínez Franco (2021 - 2022)"
$ws.Range("E12").Value = "Mariana Saavedra Botero (2021 - 2022)"
$ws.Range("E13").Value = "John Jairo Rubio (2021 - 2022)"

# Rows 14-26 - remaining "why" entries, shifted down from old rows 7-19
$ws.Range("E14").Value = "Maria Paula Moreno Rodríguez (2019 - 2021)"
$ws.Range("E15").Value = "Andrés Felipe Orozco Serrato (2020 - 2021)"
$ws.Range("E16").Value = "Danny Ferley Gaitan Rodríguez (2019 - 2020)"
$ws.Range("E17").Value = "Hasbleidy Gamboa Ordoñez (2019 - 2020)"
$ws.Range("E18").Value = "Paula Andrea Betancourt Velandia  (2018 - 2019)"
$ws.Range("B19").Value = " "
$ws.Range("E19").Value = "Ana Sofía Gómez Castelblanco (2018 - 2019)"
$ws.Range("E20").Value = "Lina María García Hoyos  (2016 - 2017)"
$ws.Range("E21").Value = "Angie Liliana Pérez Rodríguez  (2016 - 2018)"
$ws.Range("E22").Value = "Lina María Morales Sánchez (2016 - 2017)"
$ws.Range("E23").Value = "Laura Milena Estupiñan Aldana  (2016 - 2017)"
$ws.Range("E24").Value = "Vanesa Díaz Güiza  (2016 - 2018)"
$ws.Range("E25").Value = "Cindy Paola Moncada Gómez (2016 - 2017)"
$ws.Range("E26").Value = "Haydn Ricardo Roldán Morales (2015 - 2016)"

# Update the visible selection to match the reorganised view.
$ws.Range("A6:XFD8").Select() | Out-Null
